$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(6, 10).Value = 0.5411954939918103
$ws.Cells.Item(6, 12).Value = 0.08199206275836471
$ws.Cells.Item(6, 13).Value = -0.2553208893920031
$ws.Cells.Item(6, 14).Value = 0.2553208893920031

$ws.Cells.Item(7, 8).Value = 0.5781195834465759
$ws.Cells.Item(7, 9).Value = 0.1257466917056539
$ws.Cells.Item(7, 10).Value = 0.2961337248477702
$ws.Cells.Item(7, 11).Value = 0.253888855023847
$ws.Cells.Item(7, 12).Value = -0.008827085879807017
$ws.Cells.Item(7, 13).Value = -0.2450617691440402
$ws.Cells.Item(7, 14).Value = 0.253888855023847
$ws.Cells.Item(7, 15).Value = 0.5077777100476941
$ws.Cells.Item(7, 16).Value = 0.0229073098942692

$ws.Cells.Item(10, 9).Value = 0.1413896364471977
$ws.Cells.Item(10, 10).Value = 0.3398468476371659
$ws.Cells.Item(10, 12).Value = 0.0224469051696646
$ws.Cells.Item(10, 13).Value = -0.2582297336798002
$ws.Cells.Item(10, 14).Value = 0.2582297336798002

$ws.Cells.Item(11, 9).Value = 0.06103286384976526
$ws.Cells.Item(11, 12).Value = -0.08035677259743239
$ws.Cells.Item(11, 13).Value = -0.1760214836570548
$ws.Cells.Item(11, 15).Value = 0.5127565125089745

$ws.Cells.Item(12, 10).Value = 0.06148312146823425
$ws.Cells.Item(12, 12).Value = -0.06103286384976526

$ws.Cells.Item(13, 9).Value = 0.06103286384976526
$ws.Cells.Item(13, 12).Value = 0.06103286384976526

$ws.Cells.Item(17, 8).Value = 0.365712976269027
$ws.Cells.Item(17, 9).Value = 0.1486661016051938
$ws.Cells.Item(17, 10).Value = 0.4856209221257792
$ws.Cells.Item(17, 11).Value = -0.2581720544644431
$ws.Cells.Item(17, 12).Value = 0.03716010427067583
$ws.Cells.Item(17, 13).Value = 0.2210119501937673
$ws.Cells.Item(17, 14).Value = 0.2581720544644431
$ws.Cells.Item(17, 15).Value = 0.5163441089288863
$ws.Cells.Item(17, 16).Value = 0.006759163657951913

$ws.Cells.Item(19, 8).Value = 0.7479438093251026
$ws.Cells.Item(19, 9).Value = 0.07083375926355839
$ws.Cells.Item(19, 10).Value = 0.181222431411339
$ws.Cells.Item(19, 11).Value = 0.264405210784569
$ws.Cells.Item(19, 12).Value = -0.07882610468201984
$ws.Cells.Item(19, 13).Value = -0.185579106102549
$ws.Cells.Item(19, 14).Value = 0.264405210784569
$ws.Cells.Item(19, 15).Value = 0.5288104215691378
$ws.Cells.Item(19, 16).Value = 0.007457797536929802

$ws.Cells.Item(20, 8).Value = 0.9314320205147569
$ws.Cells.Item(20, 10).Value = 0.06856797948524314
$ws.Cells.Item(20, 11).Value = 0.1834882111896543
$ws.Cells.Item(20, 12).Value = -0.07083375926355839
$ws.Cells.Item(20, 13).Value = -0.1126544519260958
$ws.Cells.Item(20, 14).Value = 0.1834882111896543
$ws.Cells.Item(20, 15).Value = 0.3669764223793085
$ws.Cells.Item(20, 16).Value = 0.002365422654613125

$ws.Cells.Item(28, 8).Value = 0.1035582536879293
$ws.Cells.Item(28, 9).Value = 0.01101037230328506
$ws.Cells.Item(28, 10).Value = 0.8854313740087856
$ws.Cells.Item(28, 11).Value = -0.1341758200281135
$ws.Cells.Item(28, 12).Value = -0.08889961869581485
$ws.Cells.Item(28, 13).Value = 0.2230754387239283
$ws.Cells.Item(28, 14).Value = 0.2230754387239283
$ws.Cells.Item(28, 15).Value = 0.4461508774478567
$ws.Cells.Item(28, 16).Value = 0.004559238704794746

$ws.Cells.Item(42, 8).Value = 0.553073778772501
$ws.Cells.Item(42, 9).Value = 0.1309706264225086
$ws.Cells.Item(42, 10).Value = 0.3159555948049905
$ws.Cells.Item(42, 11).Value = -0.25442842847425
$ws.Cells.Item(42, 12).Value = 0.0828982799208427
$ws.Cells.Item(42, 13).Value = 0.1715301485534071
$ws.Cells.Item(42, 14).Value = 0.25442842847425
$ws.Cells.Item(42, 15).Value = 0.5088568569484998
$ws.Cells.Item(42, 16).Value = 0.01710009340326598

$ws.Cells.Item(43, 8).Value = 0.8088580434252878
$ws.Cells.Item(43, 9).Value = 0.04791349487671182
$ws.Cells.Item(43, 10).Value = 0.1432284616980002
$ws.Cells.Item(43, 11).Value = 0.2557842646527868
$ws.Cells.Item(43, 12).Value = -0.08305713154579675
$ws.Cells.Item(43, 13).Value = -0.1727271331069902
$ws.Cells.Item(43, 14).Value = 0.2557842646527868
$ws.Cells.Item(43, 15).Value = 0.5115685293055738
$ws.Cells.Item(43, 16).Value = 0.009568550232727078

$ws.Cells.Item(49, 8).Value = 0.2389536299135206
$ws.Cells.Item(49, 9).Value = 0.1021100533848839
$ws.Cells.Item(49, 10).Value = 0.6589363167015954
$ws.Cells.Item(49, 11).Value = -0.2179318215540149
$ws.Cells.Item(49, 12).Value = -0.05436400819250961
$ws.Cells.Item(49, 13).Value = 0.2722958297465245
$ws.Cells.Item(49, 14).Value = 0.2722958297465245
$ws.Cells.Item(49, 15).Value = 0.5445916594930491
$ws.Cells.Item(49, 16).Value = 0.01291830893320134

$ws.Cells.Item(50, 11).Value = -0.1329180762647085
$ws.Cells.Item(50, 12).Value = -0.08770532053474728
$ws.Cells.Item(50, 13).Value = 0.220623396799456
$ws.Cells.Item(50, 14).Value = 0.220623396799456
$ws.Cells.Item(50, 15).Value = 0.4412467935989118

$ws.Cells.Item(63, 8).Value = 0.3130742614062856
$ws.Cells.Item(63, 9).Value = 0.131033838861021
$ws.Cells.Item(63, 10).Value = 0.5558918997326934
$ws.Cells.Item(63, 11).Value = 0.1680112014594255
$ws.Cells.Item(63, 12).Value = 0.08250862477919413
$ws.Cells.Item(63, 13).Value = -0.2505198262386196
$ws.Cells.Item(63, 14).Value = 0.2505198262386196
$ws.Cells.Item(63, 15).Value = 0.5010396524772393
$ws.Cells.Item(63, 16).Value = 0.01706070624384681

$ws.Cells.Item(64, 11).Value = -0.1727682860510911
$ws.Cells.Item(64, 12).Value = -0.08314208769414984
$ws.Cells.Item(64, 13).Value = 0.2559103737452408
$ws.Cells.Item(64, 14).Value = 0.2559103737452408
$ws.Cells.Item(64, 15).Value = 0.5118207474904817

$ws.Cells.Item(66, 8).Value = 0.2396446960338722
$ws.Cells.Item(66, 9).Value = 0.101316779340817
$ws.Cells.Item(66, 10).Value = 0.6590385246253107
$ws.Cells.Item(66, 11).Value = -0.213048915732339
$ws.Cells.Item(66, 12).Value = -0.05622323750838265
$ws.Cells.Item(66, 13).Value = 0.2692721532407216
$ws.Cells.Item(66, 14).Value = 0.2692721532407216
$ws.Cells.Item(66, 15).Value = 0.5385443064814432
$ws.Cells.Item(66, 16).Value = 0.01936563817014287

$ws.Cells.Item(67, 8).Value = 0.1018119831948096
$ws.Cells.Item(67, 9).Value = 0.01278320438915138
$ws.Cells.Item(67, 10).Value = 0.8854048124160391
$ws.Cells.Item(67, 11).Value = -0.1378327128390626
$ws.Cells.Item(67, 12).Value = -0.08853357495166567
$ws.Cells.Item(67, 13).Value = 0.2263662877907284
$ws.Cells.Item(67, 14).Value = 0.2263662877907284
$ws.Cells.Item(67, 15).Value = 0.4527325755814567
$ws.Cells.Item(67, 16).Value = 0.00799280330776681

$ws.Cells.Item(68, 8).Value = 0.2396446960338722
$ws.Cells.Item(68, 9).Value = 0.101316779340817
$ws.Cells.Item(68, 10).Value = 0.6590385246253107
$ws.Cells.Item(68, 11).Value = 0.1378327128390626
$ws.Cells.Item(68, 12).Value = 0.08853357495166567
$ws.Cells.Item(68, 13).Value = -0.2263662877907284
$ws.Cells.Item(68, 14).Value = 0.2263662877907284
$ws.Cells.Item(68, 15).Value = 0.4527325755814567
$ws.Cells.Item(68, 16).Value = 0.01936563817014287

$ws.Cells.Item(69, 8).Value = 0.4505719099418108
$ws.Cells.Item(69, 9).Value = 0.1575896002424802
$ws.Cells.Item(69, 10).Value = 0.391838489815709
$ws.Cells.Item(69, 11).Value = 0.2109272139079386
$ws.Cells.Item(69, 12).Value = 0.05627282090166319
$ws.Cells.Item(69, 13).Value = -0.2672000348096016
$ws.Cells.Item(69, 14).Value = 0.2672000348096016
$ws.Cells.Item(69, 15).Value = 0.5344000696192034
$ws.Cells.Item(69, 16).Value = 0.02405717493533467

$ws.Cells.Item(70, 11).Value = 0.2725160406397507
$ws.Cells.Item(70, 12).Value = -0.07755280171442135
$ws.Cells.Item(70, 13).Value = -0.1949632389253295
$ws.Cells.Item(70, 14).Value = 0.2725160406397507
$ws.Cells.Item(70, 15).Value = 0.5450320812795016

$ws.Cells.Item(74, 8).Value = 0.1072569224041067
$ws.Cells.Item(74, 9).Value = 0.01986187109616585
$ws.Cells.Item(74, 10).Value = 0.8728812064997274
$ws.Cells.Item(74, 11).Value = -0.1479158552974995
$ws.Cells.Item(74, 12).Value = -0.08767940511980918
$ws.Cells.Item(74, 13).Value = 0.2355952604173086
$ws.Cells.Item(74, 14).Value = 0.2355952604173086
$ws.Cells.Item(74, 15).Value = 0.4711905208346172
$ws.Cells.Item(74, 16).Value = 0.005863922244769527

$ws.Cells.Item(75, 8).Value = 0.03668992734340081
$ws.Cells.Item(75, 10).Value = 0.9633100726565991
$ws.Cells.Item(75, 11).Value = -0.07056699506070593
$ws.Cells.Item(75, 12).Value = -0.01986187109616585
$ws.Cells.Item(75, 13).Value = 0.09042886615687162
$ws.Cells.Item(75, 14).Value = 0.09042886615687162
$ws.Cells.Item(75, 15).Value = 0.1808577323137434
$ws.Cells.Item(75, 16).Value = 0.001832640266848577

$ws.Cells.Item(77, 8).Value = 0.08547742840057443
$ws.Cells.Item(77, 10).Value = 0.9145225715994256
$ws.Cells.Item(77, 11).Value = -0.1244799445117707
$ws.Cells.Item(77, 13).Value = 0.2112379353793506
$ws.Cells.Item(77, 14).Value = 0.2112379353793506
$ws.Cells.Item(77, 15).Value = 0.4224758707587013
$ws.Cells.Item(77, 16).Value = 0.004053313062534767

$ws.Cells.Item(78, 11).Value = -0.05592070187443837
$ws.Cells.Item(78, 13).Value = 0.05592070187443832
$ws.Cells.Item(78, 14).Value = 0.05592070187443837
$ws.Cells.Item(78, 15).Value = 0.1118414037488767

$ws.Cells.Item(81, 8).Value = 0.2721721121744075
$ws.Cells.Item(81, 9).Value = 0.1140649744918982
$ws.Cells.Item(81, 10).Value = 0.6137629133336943
$ws.Cells.Item(81, 11).Value = -0.2240306009019178
$ws.Cells.Item(81, 12).Value = -0.03232853809879747
$ws.Cells.Item(81, 13).Value = 0.2563591390007151
$ws.Cells.Item(81, 14).Value = 0.2563591390007151
$ws.Cells.Item(81, 15).Value = 0.5127182780014303
$ws.Cells.Item(81, 16).Value = 0.01787355363882735

$ws.Cells.Item(82, 8).Value = 0.4983378827953688
$ws.Cells.Item(82, 9).Value = 0.1462919867281731
$ws.Cells.Item(82, 10).Value = 0.3553701304764582
$ws.Cells.Item(82, 11).Value = 0.2261657706209613
$ws.Cells.Item(82, 12).Value = 0.03222701223627487
$ws.Cells.Item(82, 13).Value = -0.2583927828572361
$ws.Cells.Item(82, 14).Value = 0.2583927828572361
$ws.Cells.Item(82, 15).Value = 0.5167855657144723
$ws.Cells.Item(82, 16).Value = 0.02013233598517662

$ws.Cells.Item(83, 8).Value = 0.7678689391106691
$ws.Cells.Item(83, 9).Value = 0.06576380598993567
$ws.Cells.Item(83, 10).Value = 0.1663672548993952
$ws.Cells.Item(83, 11).Value = 0.2695310563153003
$ws.Cells.Item(83, 12).Value = -0.08052818073823742
$ws.Cells.Item(83, 13).Value = -0.189002875577063
$ws.Cells.Item(83, 14).Value = 0.2695310563153003
$ws.Cells.Item(83, 15).Value = 0.5390621126306008
$ws.Cells.Item(83, 16).Value = 0.01261247835560044

$ws.Cells.Item(84, 8).Value = 0.506849028222723
$ws.Cells.Item(84, 9).Value = 0.1458872875734123
$ws.Cells.Item(84, 10).Value = 0.3472636842038647
$ws.Cells.Item(84, 11).Value = -0.2610199108879461
$ws.Cells.Item(84, 12).Value = 0.08012348158347665
$ws.Cells.Item(84, 13).Value = 0.1808964293044695
$ws.Cells.Item(84, 14).Value = 0.2610199108879461
$ws.Cells.Item(84, 15).Value = 0.5220398217758921
$ws.Cells.Item(84, 16).Value = 0.02004096318491041

$ws.Cells.Item(85, 11).Value = -0.2217453895642903
$ws.Cells.Item(85, 12).Value = -0.03042393416311646
$ws.Cells.Item(85, 13).Value = 0.2521693237274067
$ws.Cells.Item(85, 14).Value = 0.2521693237274067
$ws.Cells.Item(85, 15).Value = 0.5043386474548135

$ws.Cells.Item(92, 10).Value = 0.3960079295022954
$ws.Cells.Item(92, 12).Value = 0.07709505820676898
$ws.Cells.Item(92, 15).Value = 0.5427740561395203

$ws.Cells.Item(93, 8).Value = 0.7157152330212463
$ws.Cells.Item(93, 11).Value = 0.2713870280697601
$ws.Cells.Item(93, 12).Value = -0.07709505820676898
$ws.Cells.Item(93, 14).Value = 0.2713870280697601
$ws.Cells.Item(93, 16).Value = 0.009783217086726193

$ws.Cells.Item(95, 8).Value = 0.5670389009009226
$ws.Cells.Item(95, 9).Value = 0.12881444559789
$ws.Cells.Item(95, 10).Value = 0.3041466535011876
$ws.Cells.Item(95, 11).Value = 0.2509962605951318
$ws.Cells.Item(95, 12).Value = -0.002753075027381396
$ws.Cells.Item(95, 13).Value = -0.2482431855677502
$ws.Cells.Item(95, 14).Value = 0.2509962605951318
$ws.Cells.Item(95, 15).Value = 0.5019925211902634
$ws.Cells.Item(95, 16).Value = 0.01265263414743135

$ws.Cells.Item(96, 11).Value = 0.2448073533908302
$ws.Cells.Item(96, 12).Value = -0.08255788861171925
$ws.Cells.Item(96, 13).Value = -0.1622494647791111
$ws.Cells.Item(96, 14).Value = 0.2448073533908302
$ws.Cells.Item(96, 15).Value = 0.4896147067816605

$ws.Cells.Item(112, 8).Value = 0.5558669876920941
$ws.Cells.Item(112, 9).Value = 0.1302167147096241
$ws.Cells.Item(112, 10).Value = 0.3139162975982819
$ws.Cells.Item(112, 11).Value = -0.2538108130156201
$ws.Cells.Item(112, 12).Value = 0.08305139789209193
$ws.Cells.Item(112, 13).Value = 0.1707594151235284
$ws.Cells.Item(112, 14).Value = 0.2538108130156201
$ws.Cells.Item(112, 15).Value = 0.5076216260312405
$ws.Cells.Item(112, 16).Value = 0.01705220910537886

$ws.Cells.Item(113, 8).Value = 0.31205823299909
$ws.Cells.Item(113, 9).Value = 0.1306551342960267
$ws.Cells.Item(113, 10).Value = 0.5572866327048832
$ws.Cells.Item(113, 11).Value = -0.2438087546930041
$ws.Cells.Item(113, 12).Value = 0.0004384195864026441
$ws.Cells.Item(113, 13).Value = 0.2433703351066013
$ws.Cells.Item(113, 14).Value = 0.2438087546930041
$ws.Cells.Item(113, 15).Value = 0.4876175093860081
$ws.Cells.Item(113, 16).Value = 0.01703645938097853

$ws.Cells.Item(116, 8).Value = 0.8902405805161542
$ws.Cells.Item(116, 9).Value = 0.009448001818009746
$ws.Cells.Item(116, 10).Value = 0.1003114176658361
$ws.Cells.Item(116, 11).Value = 0.2239762685371836
$ws.Cells.Item(116, 12).Value = -0.08924510318424353
$ws.Cells.Item(116, 13).Value = -0.13473116535294
$ws.Cells.Item(116, 14).Value = 0.2239762685371836
$ws.Cells.Item(116, 15).Value = 0.4479525370743671
$ws.Cells.Item(116, 16).Value = 0.003654075250921813

$ws.Cells.Item(117, 8).Value = 0.6527975785770321
$ws.Cells.Item(117, 9).Value = 0.1002054823495048
$ws.Cells.Item(117, 10).Value = 0.246996939073463
$ws.Cells.Item(117, 11).Value = -0.237443001939122
$ws.Cells.Item(117, 12).Value = 0.09075748053149506
$ws.Cells.Item(117, 13).Value = 0.1466855214076269
$ws.Cells.Item(117, 14).Value = 0.237443001939122
$ws.Cells.Item(117, 15).Value = 0.474886003878244
$ws.Cells.Item(117, 16).Value = 0.009311235088877899

$ws.Cells.Item(123, 16).Value = 0.01325599429911546

$ws.Cells.Item(124, 8).Value = 0.05551547937090666
$ws.Cells.Item(124, 10).Value = 0.9444845206290934
$ws.Cells.Item(124, 11).Value = -0.09004315064595005
$ws.Cells.Item(124, 13).Value = 0.1375249435455098
$ws.Cells.Item(124, 14).Value = 0.1375249435455098
$ws.Cells.Item(124, 15).Value = 0.2750498870910196
$ws.Cells.Item(124, 16).Value = 0.00427236014912871

$ws.Cells.Item(125, 8).Value = 0.1443537408265577
$ws.Cells.Item(125, 9).Value = 0.04732302575041582
$ws.Cells.Item(125, 10).Value = 0.8083232334230265
$ws.Cells.Item(125, 11).Value = 0.088838261455651
$ws.Cells.Item(125, 12).Value = 0.04732302575041582
$ws.Cells.Item(125, 13).Value = -0.1361612872060669
$ws.Cells.Item(125, 14).Value = 0.1361612872060669
$ws.Cells.Item(125, 15).Value = 0.2723225744121337
$ws.Cells.Item(125, 16).Value = 0.01318109951692307

$ws.Cells.Item(126, 8).Value = 0.05345452382807847
$ws.Cells.Item(126, 10).Value = 0.9465454761719215
$ws.Cells.Item(126, 11).Value = -0.09089921699847919
$ws.Cells.Item(126, 12).Value = -0.04732302575041582
$ws.Cells.Item(126, 13).Value = 0.1382222427488951
$ws.Cells.Item(126, 14).Value = 0.1382222427488951
$ws.Cells.Item(126, 15).Value = 0.27644448549779
$ws.Cells.Item(126, 16).Value = 0.004122729739365263

$ws.Cells.Item(134, 8).Value = 0.5222935549413155
$ws.Cells.Item(134, 9).Value = 0.1408965399874991
$ws.Cells.Item(134, 10).Value = 0.3368099050711855
$ws.Cells.Item(134, 11).Value = -0.2539897272982634
$ws.Cells.Item(134, 12).Value = 0.08030471283860288
$ws.Cells.Item(134, 13).Value = 0.1736850144596607
$ws.Cells.Item(134, 14).Value = 0.2539897272982634
$ws.Cells.Item(134, 15).Value = 0.5079794545965269
$ws.Cells.Item(134, 16).Value = 0.01759753171354824

$ws.Cells.Item(135, 8).Value = 0.7894338698532433
$ws.Cells.Item(135, 9).Value = 0.05896486117409268
$ws.Cells.Item(135, 10).Value = 0.151601268972664
$ws.Cells.Item(135, 11).Value = 0.2671403149119278
$ws.Cells.Item(135, 12).Value = -0.08193167881340639
$ws.Cells.Item(135, 13).Value = -0.1852086360985216
$ws.Cells.Item(135, 14).Value = 0.2671403149119278
$ws.Cells.Item(135, 15).Value = 0.5342806298238558
$ws.Cells.Item(135, 16).Value = 0.01038027749692238

$ws.Cells.Item(136, 8).Value = 0.5265303236718624
$ws.Cells.Item(136, 9).Value = 0.1406537211124254
$ws.Cells.Item(136, 10).Value = 0.3328159552157123
$ws.Cells.Item(136, 11).Value = -0.2629035461813809
$ws.Cells.Item(136, 12).Value = 0.08168885993833273
$ws.Cells.Item(136, 13).Value = 0.1812146862430483
$ws.Cells.Item(136, 14).Value = 0.2629035461813809
$ws.Cells.Item(136, 15).Value = 0.525807092362762
$ws.Cells.Item(136, 16).Value = 0.01754713745022474

$ws.Cells.Item(148, 8).Value = 0.990123837775249
$ws.Cells.Item(148, 10).Value = 0.009876162224750927
$ws.Cells.Item(148, 11).Value = 0.0001009172245038936
$ws.Cells.Item(148, 13).Value = -0.0001009172245039092
$ws.Cells.Item(148, 14).Value = 0.0001009172245039092
$ws.Cells.Item(148, 15).Value = 0.0002018344490078029
$ws.Cells.Item(148, 16).Value = 0.0005070397445276254

$ws.Cells.Item(153, 8).Value = 0.1347241067216041
$ws.Cells.Item(153, 9).Value = 0.04308033841490821
$ws.Cells.Item(153, 10).Value = 0.8221955548634877
$ws.Cells.Item(153, 11).Value = -0.1691501491713128
$ws.Cells.Item(153, 12).Value = -0.08393843068155624
$ws.Cells.Item(153, 13).Value = 0.253088579852869
$ws.Cells.Item(153, 14).Value = 0.253088579852869
$ws.Cells.Item(153, 15).Value = 0.5061771597057381
$ws.Cells.Item(153, 16).Value = 0.009007050935491687

$ws.Cells.Item(154, 8).Value = 0.0492640491396921
$ws.Cells.Item(154, 10).Value = 0.9507359508603079
$ws.Cells.Item(154, 11).Value = -0.08546005758191202
$ws.Cells.Item(154, 12).Value = -0.04308033841490821
$ws.Cells.Item(154, 13).Value = 0.1285403959968202
$ws.Cells.Item(154, 14).Value = 0.1285403959968202
$ws.Cells.Item(154, 15).Value = 0.2570807919936404
$ws.Cells.Item(154, 16).Value = 0.00277553200604765

$ws.Cells.Item(155, 8).Value = 0.1347241067216041
$ws.Cells.Item(155, 9).Value = 0.04308033841490821
$ws.Cells.Item(155, 10).Value = 0.8221955548634877
$ws.Cells.Item(155, 11).Value = 0.08546005758191202
$ws.Cells.Item(155, 12).Value = 0.04308033841490821
$ws.Cells.Item(155, 13).Value = -0.1285403959968202
$ws.Cells.Item(155, 14).Value = 0.1285403959968202
$ws.Cells.Item(155, 15).Value = 0.2570807919936404
$ws.Cells.Item(155, 16).Value = 0.009007050935491687

$ws.Cells.Item(156, 11).Value = 0.1577070521518429
$ws.Cells.Item(156, 12).Value = 0.08293572320130704
$ws.Cells.Item(156, 13).Value = -0.2406427753531499
$ws.Cells.Item(156, 14).Value = 0.2406427753531499
$ws.Cells.Item(156, 15).Value = 0.4812855507062999

$ws.Cells.Item(159, 8).Value = 0.03628518049041867
$ws.Cells.Item(159, 10).Value = 0.9637148195095814
$ws.Cells.Item(159, 11).Value = -0.0728323814343696
$ws.Cells.Item(159, 13).Value = 0.09195940448488471
$ws.Cells.Item(159, 14).Value = 0.09195940448488471
$ws.Cells.Item(159, 15).Value = 0.1839188089697693
$ws.Cells.Item(159, 16).Value = 0.0007770792481599198

$ws.Cells.Item(163, 8).Value = 0.186548863766279
$ws.Cells.Item(163, 9).Value = 0.07858938293049667
$ws.Cells.Item(163, 10).Value = 0.7348617533032242
$ws.Cells.Item(163, 11).Value = -0.20321750761831
$ws.Cells.Item(163, 12).Value = -0.07895063391870305
$ws.Cells.Item(163, 13).Value = 0.2821681415370129
$ws.Cells.Item(163, 14).Value = 0.2821681415370129
$ws.Cells.Item(163, 15).Value = 0.5643362830740259
$ws.Cells.Item(163, 16).Value = 0.007759285812041206

$ws.Cells.Item(164, 11).Value = -0.1101461877854348
$ws.Cells.Item(164, 12).Value = -0.07858938293049667
$ws.Cells.Item(164, 13).Value = 0.1887355707159316
$ws.Cells.Item(164, 14).Value = 0.1887355707159316
$ws.Cells.Item(164, 15).Value = 0.3774711414318631

$ws.Cells.Item(167, 8).Value = 0.9348521692958586
$ws.Cells.Item(167, 10).Value = 0.06514783070414154
$ws.Cells.Item(167, 11).Value = 0.1700671005178169
$ws.Cells.Item(167, 13).Value = -0.1050834633041654
$ws.Cells.Item(167, 14).Value = 0.1700671005178169
$ws.Cells.Item(167, 15).Value = 0.3401342010356336
$ws.Cells.Item(167, 16).Value = 0.002255688550321704

$ws.Cells.Item(173, 8).Value = 0.7137523445297103
$ws.Cells.Item(173, 9).Value = 0.08361727176794841
$ws.Cells.Item(173, 10).Value = 0.2026303837023413
$ws.Cells.Item(173, 11).Value = 0.2735966910503468
$ws.Cells.Item(173, 12).Value = -0.07710450750187954
$ws.Cells.Item(173, 13).Value = -0.1964921835484675
$ws.Cells.Item(173, 14).Value = 0.2735966910503468
$ws.Cells.Item(173, 15).Value = 0.5471933821006938
$ws.Cells.Item(173, 16).Value = 0.009833481558690417

$ws.Cells.Item(174, 8).Value = 0.9191197839661166
$ws.Cells.Item(174, 10).Value = 0.08088021603388329
$ws.Cells.Item(174, 11).Value = 0.2053674394364063
$ws.Cells.Item(174, 12).Value = -0.08361727176794841
$ws.Cells.Item(174, 13).Value = -0.121750167668458
$ws.Cells.Item(174, 14).Value = 0.2053674394364063
$ws.Cells.Item(174, 15).Value = 0.4107348788728127
$ws.Cells.Item(174, 16).Value = 0.00330393807503092

$ws.Cells.Item(178, 13).Value = -0.2277897700998084
$ws.Cells.Item(178, 15).Value = 0.5002039435734377

$ws.Cells.Item(182, 8).Value = 0.07782632795625791
$ws.Cells.Item(182, 11).Value = -0.1190489229341216
$ws.Cells.Item(182, 15).Value = 0.3981714429243611

$ws.Cells.Item(186, 8).Value = 0.4326173638030108
$ws.Cells.Item(186, 9).Value = 0.1652357624588611
$ws.Cells.Item(186, 10).Value = 0.402146873738128
$ws.Cells.Item(186, 11).Value = 0.2107651696752403
$ws.Cells.Item(186, 12).Value = 0.07267677879279942
$ws.Cells.Item(186, 13).Value = -0.2834419484680397
$ws.Cells.Item(186, 14).Value = 0.2834419484680397
$ws.Cells.Item(186, 15).Value = 0.5668838969360794
$ws.Cells.Item(186, 16).Value = 0.01617303984808872

$ws.Cells.Item(188, 8).Value = 0.06639471196549736
$ws.Cells.Item(188, 10).Value = 0.9336052880345026
$ws.Cells.Item(188, 11).Value = -0.1089454573413186
$ws.Cells.Item(188, 13).Value = 0.17697900813349
$ws.Cells.Item(188, 14).Value = 0.17697900813349
$ws.Cells.Item(188, 15).Value = 0.35395801626698
$ws.Cells.Item(188, 16).Value = 0.001836635679659734

$ws.Cells.Item(195, 8).Value = 0.7228440979870773
$ws.Cells.Item(195, 11).Value = 0.2743364043527137
$ws.Cells.Item(195, 14).Value = 0.2743364043527137
$ws.Cells.Item(195, 15).Value = 0.5486728087054276
$ws.Cells.Item(195, 16).Value = 0.006406392514215136

$ws.Cells.Item(196, 11).Value = 0.1951186270928609
$ws.Cells.Item(196, 14).Value = 0.1951186270928609
$ws.Cells.Item(196, 15).Value = 0.3902372541857217

$ws.Cells.Item(198, 8).Value = 0.3594402428258734
$ws.Cells.Item(198, 9).Value = 0.1464951794679015
$ws.Cells.Item(198, 10).Value = 0.4940645777062251
$ws.Cells.Item(198, 11).Value = -0.256341514134253
$ws.Cells.Item(198, 12).Value = 0.03262721314091876
$ws.Cells.Item(198, 13).Value = 0.2237143009933342
$ws.Cells.Item(198, 14).Value = 0.256341514134253
$ws.Cells.Item(198, 15).Value = 0.512683028268506
$ws.Cells.Item(198, 16).Value = 0.006724911858732378

$ws.Cells.Item(201, 8).Value = 0.1482104339992398
$ws.Cells.Item(201, 9).Value = 0.05004541922573329
$ws.Cells.Item(201, 10).Value = 0.801744146775027
$ws.Cells.Item(201, 11).Value = -0.1698739911016994
$ws.Cells.Item(201, 12).Value = -0.08227562198901724
$ws.Cells.Item(201, 13).Value = 0.2521496130907167
$ws.Cells.Item(201, 14).Value = 0.2521496130907167
$ws.Cells.Item(201, 15).Value = 0.5042992261814333
$ws.Cells.Item(201, 16).Value = 0.00739412103067248

$ws.Cells.Item(202, 11).Value = -0.09139445053734992
$ws.Cells.Item(202, 12).Value = -0.05004541922573329
$ws.Cells.Item(202, 13).Value = 0.1414398697630831
$ws.Cells.Item(202, 14).Value = 0.1414398697630831
$ws.Cells.Item(202, 15).Value = 0.2828797395261663

$ws.Cells.Item(206, 8).Value = 0.2747155454308867
$ws.Cells.Item(206, 9).Value = 0.1184049296349224
$ws.Cells.Item(206, 10).Value = 0.606879524934191
$ws.Cells.Item(206, 11).Value = -0.2411953570079378
$ws.Cells.Item(206, 12).Value = -0.02285741471887752
$ws.Cells.Item(206, 13).Value = 0.2640527717268152
$ws.Cells.Item(206, 14).Value = 0.2640527717268152
$ws.Cells.Item(206, 15).Value = 0.5281055434536306
$ws.Cells.Item(206, 16).Value = 0.00602454315502678

$ws.Cells.Item(209, 8).Value = 0.2387417121874942
$ws.Cells.Item(209, 9).Value = 0.1009126927039771
$ws.Cells.Item(209, 10).Value = 0.6603455951085286
$ws.Cells.Item(209, 11).Value = 0.1364962878894062
$ws.Cells.Item(209, 12).Value = 0.08856701369163142
$ws.Cells.Item(209, 13).Value = -0.2250633015810376
$ws.Cells.Item(209, 14).Value = 0.2250633015810376
$ws.Cells.Item(209, 15).Value = 0.4501266031620752
$ws.Cells.Item(209, 16).Value = 0.01103917151852843

$ws.Cells.Item(210, 8).Value = 0.4321656057108492
$ws.Cells.Item(210, 9).Value = 0.1583257072067548
$ws.Cells.Item(210, 10).Value = 0.4095086870823961
$ws.Cells.Item(210, 11).Value = 0.1934238935233549
$ws.Cells.Item(210, 12).Value = 0.05741301450277767
$ws.Cells.Item(210, 13).Value = -0.2508369080261326
$ws.Cells.Item(210, 14).Value = 0.2508369080261326
$ws.Cells.Item(210, 15).Value = 0.5016738160522651
$ws.Cells.Item(210, 16).Value = 0.01378818877515796

$ws.Cells.Item(219, 8).Value = 0.4011886883672275
$ws.Cells.Item(219, 9).Value = 0.1621027960930984
$ws.Cells.Item(219, 10).Value = 0.436708515539674
$ws.Cells.Item(219, 11).Value = 0.1937263317844816
$ws.Cells.Item(219, 12).Value = 0.07659753774833392
$ws.Cells.Item(219, 13).Value = -0.2703238695328155
$ws.Cells.Item(219, 14).Value = 0.2703238695328155
$ws.Cells.Item(219, 15).Value = 0.540647739065631
$ws.Cells.Item(219, 16).Value = 0.02303911082518953

$ws.Cells.Item(220, 8).Value = 0.1997806865859413
$ws.Cells.Item(220, 9).Value = 0.08445053444726962
$ws.Cells.Item(220, 10).Value = 0.7157687789667891
$ws.Cells.Item(220, 11).Value = -0.2014080017812862
$ws.Cells.Item(220, 12).Value = -0.07765226164582882
$ws.Cells.Item(220, 13).Value = 0.2790602634271152
$ws.Cells.Item(220, 14).Value = 0.2790602634271152
$ws.Cells.Item(220, 15).Value = 0.5581205268542302
$ws.Cells.Item(220, 16).Value = 0.01631966072426712

$ws.Cells.Item(221, 8).Value = 0.3928784076831426
$ws.Cells.Item(221, 9).Value = 0.1620012824463322
$ws.Cells.Item(221, 10).Value = 0.4451203098705252
$ws.Cells.Item(221, 11).Value = 0.1930977210972013
$ws.Cells.Item(221, 12).Value = 0.07755074799906263
$ws.Cells.Item(221, 13).Value = -0.2706484690962639
$ws.Cells.Item(221, 14).Value = 0.2706484690962639
$ws.Cells.Item(221, 15).Value = 0.5412969381925279
$ws.Cells.Item(221, 16).Value = 0.02301000188899517

$ws.Cells.Item(222, 8).Value = 0.1997806865859413
$ws.Cells.Item(222, 9).Value = 0.08445053444726962
$ws.Cells.Item(222, 10).Value = 0.7157687789667891
$ws.Cells.Item(222, 11).Value = -0.1930977210972013
$ws.Cells.Item(222, 12).Value = -0.07755074799906263
$ws.Cells.Item(222, 13).Value = 0.2706484690962639
$ws.Cells.Item(222, 14).Value = 0.2706484690962639
$ws.Cells.Item(222, 15).Value = 0.5412969381925279
$ws.Cells.Item(222, 16).Value = 0.01631966072426712

$ws.Cells.Item(223, 8).Value = 0.3928784076831426
$ws.Cells.Item(223, 9).Value = 0.1620012824463322
$ws.Cells.Item(223, 10).Value = 0.4451203098705252
$ws.Cells.Item(223, 11).Value = 0.1930977210972013
$ws.Cells.Item(223, 12).Value = 0.07755074799906263
$ws.Cells.Item(223, 13).Value = -0.2706484690962639
$ws.Cells.Item(223, 14).Value = 0.2706484690962639
$ws.Cells.Item(223, 15).Value = 0.5412969381925279
$ws.Cells.Item(223, 16).Value = 0.02301000188899517

$ws.Cells.Item(224, 8).Value = 0.193770464299839
$ws.Cells.Item(224, 9).Value = 0.08362530698018947
$ws.Cells.Item(224, 10).Value = 0.7226042287199714
$ws.Cells.Item(224, 11).Value = -0.1991079433833036
$ws.Cells.Item(224, 12).Value = -0.07837597546614278
$ws.Cells.Item(224, 13).Value = 0.2774839188494462
$ws.Cells.Item(224, 14).Value = 0.2774839188494462
$ws.Cells.Item(224, 15).Value = 0.5549678376988927
$ws.Cells.Item(224, 16).Value = 0.01604825717901885

$ws.Cells.Item(227, 8).Value = 0.2875842700549563
$ws.Cells.Item(227, 9).Value = 0.1246853172682162
$ws.Cells.Item(227, 10).Value = 0.5877304126768274
$ws.Cells.Item(227, 11).Value = -0.252209638783856
$ws.Cells.Item(227, 12).Value = -0.01026277961759699
$ws.Cells.Item(227, 13).Value = 0.262472418401453
$ws.Cells.Item(227, 14).Value = 0.262472418401453
$ws.Cells.Item(227, 15).Value = 0.5249448368029059
$ws.Cells.Item(227, 16).Value = 0.02060451189960771

$ws.Cells.Item(228, 11).Value = 0.2394075818864032
$ws.Cells.Item(228, 12).Value = 0.01113694219717136
$ws.Cells.Item(228, 13).Value = -0.2505445240835745
$ws.Cells.Item(228, 14).Value = 0.2505445240835745
$ws.Cells.Item(228, 15).Value = 0.5010890481671491

$ws.Cells.Item(262, 11).Value = -2.419626224819948E-05
$ws.Cells.Item(262, 13).Value = 2.419626224814397E-05
$ws.Cells.Item(262, 14).Value = 2.419626224819948E-05
$ws.Cells.Item(262, 15).Value = 4.839252449634346E-05

$ws.Cells.Item(263, 11).Value = -7.651776088583168E-06
$ws.Cells.Item(263, 13).Value = 7.651776088524187E-06
$ws.Cells.Item(263, 14).Value = 7.651776088583168E-06
$ws.Cells.Item(263, 15).Value = 1.530355217710735E-05

$ws.Cells.Item(283, 8).Value = 0.5334613330644971
$ws.Cells.Item(283, 9).Value = 0.1362102569618645
$ws.Cells.Item(283, 10).Value = 0.3303284099736386
$ws.Cells.Item(283, 11).Value = 0.235612750259025
$ws.Cells.Item(283, 12).Value = 0.01148377993779445
$ws.Cells.Item(283, 13).Value = -0.2470965301968193
$ws.Cells.Item(283, 14).Value = 0.2470965301968193
$ws.Cells.Item(283, 15).Value = 0.4941930603936386
$ws.Cells.Item(283, 16).Value = 0.008707391312410676

$ws.Cells.Item(284, 11).Value = -0.2506889524565951
$ws.Cells.Item(284, 12).Value = -0.01286317702975552
$ws.Cells.Item(284, 13).Value = 0.2635521294863505
$ws.Cells.Item(284, 14).Value = 0.2635521294863505
$ws.Cells.Item(284, 15).Value = 0.5271042589727011

$ws.Cells.Item(288, 8).Value = 0.7139715611744353
$ws.Cells.Item(288, 9).Value = 0.08278155563668295
$ws.Cells.Item(288, 10).Value = 0.2032468831888818
$ws.Cells.Item(288, 11).Value = 0.269643356222949
$ws.Cells.Item(288, 12).Value = -0.07688230990953555
$ws.Cells.Item(288, 13).Value = -0.1927610463134135
$ws.Cells.Item(288, 14).Value = 0.269643356222949
$ws.Cells.Item(288, 15).Value = 0.5392867124458981
$ws.Cells.Item(288, 16).Value = 0.006549370790437267

$ws.Cells.Item(289, 11).Value = -0.2865904374762533
$ws.Cells.Item(289, 12).Value = 0.07719150760875436
$ws.Cells.Item(289, 13).Value = 0.2093989298674989
$ws.Cells.Item(289, 14).Value = 0.2865904374762533
$ws.Cells.Item(289, 15).Value = 0.5731808749525066

$ws.Cells.Item(295, 8).Value = 0.583626720296482
$ws.Cells.Item(295, 9).Value = 0.1250226793207048
$ws.Cells.Item(295, 10).Value = 0.2913506003828133
$ws.Cells.Item(295, 11).Value = 0.2583687260211076
$ws.Cells.Item(295, 12).Value = -0.009925417565108396
$ws.Cells.Item(295, 13).Value = -0.2484433084559991
$ws.Cells.Item(295, 14).Value = 0.2583687260211076
$ws.Cells.Item(295, 15).Value = 0.5167374520422151
$ws.Cells.Item(295, 16).Value = 0.01655893359016351

$ws.Cells.Item(296, 8).Value = 0.3351871863717139
$ws.Cells.Item(296, 9).Value = 0.135675780074875
$ws.Cells.Item(296, 10).Value = 0.5291370335534114
$ws.Cells.Item(296, 11).Value = -0.2484395339247681
$ws.Cells.Item(296, 12).Value = 0.0106531007541702
$ws.Cells.Item(296, 13).Value = 0.2377864331705981
$ws.Cells.Item(296, 14).Value = 0.2484395339247681
$ws.Cells.Item(296, 15).Value = 0.4968790678495364
$ws.Cells.Item(296, 16).Value = 0.01745942614860991

$ws.Cells.Item(297, 8).Value = 0.583626720296482
$ws.Cells.Item(297, 9).Value = 0.1250226793207047
$ws.Cells.Item(297, 10).Value = 0.2913506003828132
$ws.Cells.Item(297, 11).Value = 0.2484395339247681
$ws.Cells.Item(297, 12).Value = -0.01065310075417022
$ws.Cells.Item(297, 13).Value = -0.2377864331705981
$ws.Cells.Item(297, 14).Value = 0.2484395339247681
$ws.Cells.Item(297, 15).Value = 0.4968790678495365
$ws.Cells.Item(297, 16).Value = 0.01655893359016351

$ws.Cells.Item(298, 8).Value = 0.3312026872607655
$ws.Cells.Item(298, 9).Value = 0.1353837670901591
$ws.Cells.Item(298, 10).Value = 0.5334135456490755
$ws.Cells.Item(298, 11).Value = -0.2524240330357165
$ws.Cells.Item(298, 12).Value = 0.01036108776945441
$ws.Cells.Item(298, 13).Value = 0.2420629452662623
$ws.Cells.Item(298, 14).Value = 0.2524240330357165
$ws.Cells.Item(298, 15).Value = 0.5048480660714332
$ws.Cells.Item(298, 16).Value = 0.01740580755193502

$ws.Cells.Item(303, 8).Value = 0.2295314400009479
$ws.Cells.Item(303, 9).Value = 0.1006612912163856
$ws.Cells.Item(303, 10).Value = 0.6698072687826665
$ws.Cells.Item(303, 11).Value = -0.2259561049600518
$ws.Cells.Item(303, 12).Value = -0.05616838837889261
$ws.Cells.Item(303, 13).Value = 0.2821244933389445
$ws.Cells.Item(303, 14).Value = 0.2821244933389445
$ws.Cells.Item(303, 15).Value = 0.5642489866778889
$ws.Cells.Item(303, 16).Value = 0.01085646322639213

$ws.Cells.Item(304, 8).Value = 0.4553188847080373
$ws.Cells.Item(304, 9).Value = 0.1568339035364801
$ws.Cells.Item(304, 10).Value = 0.3878472117554827
$ws.Cells.Item(304, 11).Value = 0.2257874447070894
$ws.Cells.Item(304, 12).Value = 0.05617261232009442
$ws.Cells.Item(304, 13).Value = -0.2819600570271838
$ws.Cells.Item(304, 14).Value = 0.2819600570271838
$ws.Cells.Item(304, 15).Value = 0.5639201140543677
$ws.Cells.Item(304, 16).Value = 0.01372583067251637

$ws.Cells.Item(305, 8).Value = 0.2383436021150399
$ws.Cells.Item(305, 9).Value = 0.1017421151954591
$ws.Cells.Item(305, 10).Value = 0.6599142826895008
$ws.Cells.Item(305, 11).Value = -0.2169752825929974
$ws.Cells.Item(305, 12).Value = -0.05509178834102092
$ws.Cells.Item(305, 13).Value = 0.2720670709340181
$ws.Cells.Item(305, 14).Value = 0.2720670709340181
$ws.Cells.Item(305, 15).Value = 0.5441341418680364
$ws.Cells.Item(305, 16).Value = 0.01105231130730629

$ws.Cells.Item(308, 8).Value = 0.5124499825576856
$ws.Cells.Item(308, 9).Value = 0.1434999312890419
$ws.Cells.Item(308, 10).Value = 0.3440500861532725
$ws.Cells.Item(308, 11).Value = 0.2353300688503645
$ws.Cells.Item(308, 12).Value = 0.02689215743745182
$ws.Cells.Item(308, 13).Value = -0.2622222262878162
$ws.Cells.Item(308, 14).Value = 0.2622222262878162
$ws.Cells.Item(308, 15).Value = 0.5244444525756325
$ws.Cells.Item(308, 16).Value = 0.01329849607365764

